$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Market Cap (column C) values for rows 2-23
$ws.Range("C2").Value = 729452307914.964
$ws.Range("C3").Value = 247274651738.1406
$ws.Range("C4").Value = 35889246967.43371
$ws.Range("C5").Value = 32841118931.31339
$ws.Range("C6").Value = 24437405375.18588
$ws.Range("C7").Value = 13445294161.80837
$ws.Range("C8").Value = 10773844725.77612
$ws.Range("C9").Value = 8984553336.742424
$ws.Range("C10").Value = 8235606022.197677
$ws.Range("C11").Value = 8035539860.119489
$ws.Range("C12").Value = 7438036117.054225
$ws.Range("C13").Value = 7192503295.905393
$ws.Range("C14").Value = 6462221783.156593
$ws.Range("C15").Value = 6071707533.534866
$ws.Range("C16").Value = 5110220791.4047
$ws.Range("C17").Value = 4795370539.839595
$ws.Range("C18").Value = 4412458485.5137
$ws.Range("C19").Value = 3792471051.095762
$ws.Range("C20").Value = 3489926144.652276
$ws.Range("C21").Value = 3322514833.25357
$ws.Range("C22").Value = 3306057057.413295
$ws.Range("C23").Value = 3007520499.669557

# Rows 24 and 25 swap: Kaspa/KAS-USD and Ethereum Classic/ETC-USD swap positions
$ws.Range("A24").Value = "Ethereum Classic"
$ws.Range("B24").Value = "ETC-USD"
$ws.Range("C24").Value = 2732922419.353358

$ws.Range("A25").Value = "Kaspa"
$ws.Range("B25").Value = "KAS-USD"
$ws.Range("C25").Value = 2730235727.597404

$ws.Range("C26").Value = 2474737191.456425
